# Translation update ("update translations 15 november"):
# Appends 11 new localization rows (keys + translations for columns
# A, B, G:Q) to Sheet1, right after the existing last row (486).
#
# The "English copy" (column B, duplicated across G:Q) is written first
# for every new row, and only afterwards is the "Key" column (A) filled
# in (identical to B except for four rows where the key carries a
# trailing period) - this mirrors the shared-string insertion order
# seen in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 487

# Row data: English copy text, and whether the Key (col A) differs from
# it by a trailing period.
$rows = @(
    @{ Text = "Contribute to keep your language on top"; KeyHasPeriod = $true },
    @{ Text = "Validate to keep your language on top"; KeyHasPeriod = $true },
    @{ Text = "Please don't use only numerics or email as username"; KeyHasPeriod = $false },
    @{ Text = "Only 1000 characters allowed"; KeyHasPeriod = $false },
    @{ Text = "Contribute to see your language on top"; KeyHasPeriod = $true },
    @{ Text = "Validate to see your language on top"; KeyHasPeriod = $true },
    @{ Text = "We feel the text you entered doesn't match the original text, are you sure about your edit"; KeyHasPeriod = $false },
    @{ Text = "404 Error"; KeyHasPeriod = $false },
    @{ Text = "Seems this page doesn't exist"; KeyHasPeriod = $false },
    @{ Text = "Visit our homepage"; KeyHasPeriod = $false },
    @{ Text = "Unspecified location"; KeyHasPeriod = $false }
)

$dataCols = @(2, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17)  # B, G, H, I, J, K, L, M, N, O, P, Q

# Pass 1: write the English copy into column B and G:Q for every new row.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $text = $rows[$i].Text
    foreach ($col in $dataCols) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.Value = $text
        $cell.Style = "Normal"
    }
}

# Pass 2: write the Key column (A) for every new row.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $text = $rows[$i].Text
    if ($rows[$i].KeyHasPeriod) {
        $text = $text + "."
    }
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Update the active selection to match the new bottom of the sheet (the
# frozen pane keeps tracking the split automatically).
$ws.Range("A501").Select()
